$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update report period / issue number labels ----
$ws.Range("A8").Value = "Volume 30   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/13/2023  Through  2/19/2023"

# ---- Helper functions ----
function Set-TextCell($ref, $text) {
    $dst = $ws.Range($ref)
    $styleDonor = $ws.Range("C14")
    $dst.Formula = "=""" + $text + """"
    $dst.Copy()
    $dst.PasteSpecial(-4163)
    $styleDonor.Copy()
    $dst.PasteSpecial(-4122)
}

function Set-NumCell($ref, $val, $styleDonorRef) {
    $dst = $ws.Range($ref)
    $styleDonor = $ws.Range($styleDonorRef)
    $styleDonor.Copy()
    $dst.PasteSpecial(-4122)
    $dst.Value = $val
}

# ---- Cells that flip between numeric and "N/A"-style text ----
Set-TextCell "F14" "0"
Set-TextCell "D15" "0"
Set-TextCell "E15" "***.*"
Set-TextCell "D17" "0"
Set-TextCell "E17" "***.*"
Set-TextCell "D22" "0"
Set-TextCell "E22" "***.*"
Set-TextCell "D26" "0"
Set-TextCell "E26" "***.*"

Set-NumCell "C16" 2 "C40"
Set-NumCell "C22" 1 "C40"
Set-NumCell "F22" 1 "C40"
Set-NumCell "D30" 2 "C40"
Set-NumCell "E30" -100 "K40"
Set-NumCell "G30" 2 "C40"
Set-NumCell "H30" -100 "K40"
Set-NumCell "J30" 2 "C40"
Set-NumCell "K30" -100 "K40"

# ---- Plain numeric value updates (style unchanged) ----
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 15
$ws.Range("K16").Value = -20
$ws.Range("L16").Value = 50
$ws.Range("M16").Value = -25
$ws.Range("N16").Value = -85.185185185185
$ws.Range("F17").Value = 6
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 12
$ws.Range("K17").Value = 20
$ws.Range("L17").Value = 20
$ws.Range("M17").Value = -7.692307692307
$ws.Range("N17").Value = 33.333333333333
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 66.666666666666
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 17
$ws.Range("K18").Value = 52.941176470588
$ws.Range("L18").Value = 44.444444444444
$ws.Range("M18").Value = 4
$ws.Range("N18").Value = -78.151260504201
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 60
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = 17.647058823529
$ws.Range("I19").Value = 92
$ws.Range("J19").Value = 89
$ws.Range("K19").Value = 3.370786516853
$ws.Range("L19").Value = 73.584905660377
$ws.Range("M19").Value = -14.814814814814
$ws.Range("N19").Value = -68.918918918918
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 300
$ws.Range("I20").Value = 13
$ws.Range("J20").Value = 6
$ws.Range("K20").Value = 116.666666666667
$ws.Range("L20").Value = 116.666666666667
$ws.Range("M20").Value = 1200
$ws.Range("N20").Value = -92.441860465116
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -26.086956521739
$ws.Range("F21").Value = 99
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = 13.793103448275
$ws.Range("H22").Value = -66.666666666666
$ws.Range("I21").Value = 157
$ws.Range("I22").Value = 6
$ws.Range("J21").Value = 141
$ws.Range("K21").Value = 11.347517730496
$ws.Range("K22").Value = 50
$ws.Range("L21").Value = 61.855670103092
$ws.Range("L22").Value = 200
$ws.Range("M21").Value = -4.268292682926
$ws.Range("M22").Value = 20
$ws.Range("N21").Value = -76.911764705882
$ws.Range("F23").Value = 4
$ws.Range("I23").Value = 5
$ws.Range("K23").Value = 400
$ws.Range("L23").Value = 400
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -48.571428571428
$ws.Range("F24").Value = 78
$ws.Range("G24").Value = 106
$ws.Range("H24").Value = -26.415094339622
$ws.Range("I24").Value = 136
$ws.Range("J24").Value = 177
$ws.Range("K24").Value = -23.163841807909
$ws.Range("L24").Value = -40.088105726872
$ws.Range("M24").Value = 5.426356589147
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 50
$ws.Range("I25").Value = 30
$ws.Range("J25").Value = 25
$ws.Range("K25").Value = 20
$ws.Range("L25").Value = 66.666666666666
$ws.Range("M25").Value = 7.142857142857
$ws.Range("L26").Value = -33.333333333333
$ws.Range("G27").Value = 1
$ws.Range("L27").Value = -66.666666666666
